# Trade #17 closed at 2026-02-17 08:19:28 - unknown UNKNOWN +0.000%
#
# Updates the workbook to reflect a new closed trade:
#  - Summary sheet: bump Total Trades (B6) and recompute Win Rate % (B9)
#  - Strategy Status sheet: bump MarketMaking row Trades (D4) and Win Rate % (G4)
#  - All Trades / MarketMaking sheets: append the new trade as row 18

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 17
$summary.Range("B9").Value = 29.41

# ---------------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 17
$status.Range("G4").Value = 29.41

# ---------------------------------------------------------------------------
# Append the new trade row (row 18) to both "All Trades" and "MarketMaking"
# sheets - they mirror each other.
# ---------------------------------------------------------------------------
$newRow = 18

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($newRow, 1).Value = 17                # A: Trade #

    # Keep date/time as literal text (not auto-converted to date serials)
    $ws.Cells.Item($newRow, 2).Value = "'2026-02-17"      # B: Date
    $ws.Cells.Item($newRow, 3).Value = "'08:19:21"        # C: Time

    $ws.Cells.Item($newRow, 4).Value = "MarketMaking"     # D: Strategy
    $ws.Cells.Item($newRow, 5).Value = "UP"               # E: Side
    $ws.Cells.Item($newRow, 6).Value = 0.03               # F: Entry Price
    $ws.Cells.Item($newRow, 7).Value = 0.03               # G: Exit Price
    $ws.Cells.Item($newRow, 8).Value = "CLOSED"           # H: Status
    $ws.Cells.Item($newRow, 9).Value = 0                  # I: P&L %
    $ws.Cells.Item($newRow, 10).Value = 0                 # J: P&L $
    $ws.Cells.Item($newRow, 11).Value = 100.01            # K: Capital After
    $ws.Cells.Item($newRow, 12).Value = 0                 # L: Entry Slippage (bps)
    $ws.Cells.Item($newRow, 13).Value = 0                 # M: Exit Slippage (bps)
    $ws.Cells.Item($newRow, 14).Value = 0.6               # N: Confidence
    $ws.Cells.Item($newRow, 15).Value = "Normal spread capture: 19600 bps"  # O: Entry Reason
    $ws.Cells.Item($newRow, 16).Value = "early_exit"      # P: Exit Reason
    $ws.Cells.Item($newRow, 17).Value = 0.13              # Q: Duration (min)
}
